$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.532.96'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.40%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.670.77'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.08%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9991'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.56'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.41%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4766'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.43%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2627'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.66%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06164'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.56%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.673.19'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06983'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.62%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.85'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.5891'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -4.18%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.374'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '75.19'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +3.44%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.0000'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.515.59'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006753'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.40'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.886.38'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.39%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.439'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.729'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.262'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '136.63'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.93%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.02'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.51%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.387'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.728'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +4.35%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '104.67'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.21%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.957'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +6.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.07841'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.52%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.642'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.51%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9991'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04263'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.88%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.99%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9543'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +3.71%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6043'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +4.19%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.587'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.97%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8858'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +7.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9999'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.23%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.863'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +3.82%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01478'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -4.64%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '96.39'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.12%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3753'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.883'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +3.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1116'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.13%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.206'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05259'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.59%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '29.85'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.371'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.92%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.29%  '
